# Applies the "Done changes in Excel and Done Code Modifications" edit:
#  - Login sheet gains a FavText column and 5 more sample rows
#  - Register/Customer shared-string indices shift automatically once the
#    now-unused test@gmail.com / Test1234 strings are overwritten below
#  - Sheet2 is renamed to InvalidLoginDetails and populated with sample data,
#    and becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Login sheet: add FavText column + more rows of sample creds
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item("Login")

$login.Range("D1").Value = "FavText"

$loginValues = @("manjunath1234", "nitin111", "nitin11212", "autoenroll1", "jayeshbulk1", "santhosh123")
$r = 2
foreach ($v in $loginValues) {
    $login.Cells.Item($r, 2).Value = $v
    $login.Cells.Item($r, 3).Value = $v
    $login.Cells.Item($r, 4).Value = $v
    $r = $r + 1
}

# keep showing the old mailto text on the (now relocated) hyperlink cell
foreach ($h in $login.Hyperlinks) {
    $h.TextToDisplay = "test@gmail.com"
}

# widen the password / new FavText columns
$login.Columns.Item(3).ColumnWidth = 16
$login.Columns.Item(4).ColumnWidth = 28.5

$login.Range("A1:C7").Select()

# ---------------------------------------------------------------------
# 2. Rename Sheet2 -> InvalidLoginDetails and fill it with sample data
# ---------------------------------------------------------------------
$invalid = $wb.Worksheets.Item("Sheet2")
$invalid.Name = "InvalidLoginDetails"

$invalid.Range("A1").Value = "TC_ID"
$invalid.Range("B1").Value = "userName"
$invalid.Range("C1").Value = "password"
$invalid.Range("A1:C1").Font.Bold = $true

$invalid.Range("A2").Value = 1
$invalid.Range("B2").Value = "pawanKalyan"
$invalid.Range("C2").Value = "pawankalyan"
$invalid.Range("B2").Style = "Hyperlink"

$invalid.Range("A3").Value = 2
$invalid.Range("B3").Value = "rancheran"
$invalid.Range("C3").Value = "chinajivi"

$invalid.Range("A4").Value = 3
$invalid.Range("A5").Value = 4
$invalid.Range("A6").Value = 5
$invalid.Range("A7").Value = 6

$invalid.Columns.Item(2).ColumnWidth = 20.833333333333332
$invalid.Columns.Item(3).ColumnWidth = 19.5

# make InvalidLoginDetails the active sheet/tab, with C8 selected
$invalid.Range("C8").Select()

Write-Host "done"
